$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3290.5386
$ws.Range("J40").Value = 4249.4
$ws.Range("L40").Value = 4249.4
$ws.Range("N40").Value = -4599.4
$ws.Range("H76").Value = 4383.3076
$ws.Range("I76").Value = 4248.5835
$ws.Range("K76").Value = 4248.5835
$ws.Range("M76").Value = -3933.5835
$ws.Range("H79").Value = 4383.3076
$ws.Range("I79").Value = 4248.5835
$ws.Range("K79").Value = 4248.5835
$ws.Range("M79").Value = -3156.5835
$ws.Range("H132").Value = 2279590
$ws.Range("I132").Value = 2316318.2
$ws.Range("K132").Value = 6948954.600000001
$ws.Range("M132").Value = -6946424.600000001
$ws.Range("H135").Value = 16076.883
$ws.Range("I135").Value = 827.5925999999999
$ws.Range("K135").Value = 7448.3334
$ws.Range("M135").Value = -4913.3334
$ws.Range("H137").Value = 22155.117
$ws.Range("J137").Value = 5143.125
$ws.Range("L137").Value = 15429.375
$ws.Range("N137").Value = -20529.375
$ws.Range("H138").Value = 4006.8667
$ws.Range("J138").Value = 4271.5625
$ws.Range("L138").Value = 12814.6875
$ws.Range("N138").Value = -23094.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26044.69
$ws.Range("I32").Value = 29513.82
$ws.Range("J32").Value = 3495.3333
$ws.Range("K32").Value = 29513.82
$ws.Range("L32").Value = 3495.3333
$ws.Range("M32").Value = -29226.82
$ws.Range("N32").Value = -4069.3333
$ws.Range("H63").Value = 4029
$ws.Range("I63").Value = 1985.6666
$ws.Range("J63").Value = 4586.273
$ws.Range("K63").Value = 1985.6666
$ws.Range("L63").Value = 4586.273
$ws.Range("M63").Value = -1299.6666
$ws.Range("N63").Value = -5958.273
$ws.Range("H66").Value = 4029
$ws.Range("I66").Value = 1985.6666
$ws.Range("J66").Value = 4586.273
$ws.Range("K66").Value = 9928.333000000001
$ws.Range("L66").Value = 22931.365
$ws.Range("M66").Value = -6496.333000000001
$ws.Range("N66").Value = -29795.365
$ws.Range("H102").Value = 12402.5
$ws.Range("I102").Value = 12505
$ws.Range("K102").Value = 12505
$ws.Range("M102").Value = -10883
$ws.Range("H122").Value = 2437.1428
$ws.Range("I122").Value = 2373.1177
$ws.Range("K122").Value = 7119.353099999999
$ws.Range("M122").Value = -4669.353099999999
$ws.Range("H132").Value = 893.5814
$ws.Range("I132").Value = 795.2432
$ws.Range("K132").Value = 2385.7296
$ws.Range("M132").Value = 144.2703999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 353
$ws.Range("J22").Value = 530.3333
$ws.Range("L22").Value = 530.3333
$ws.Range("N22").Value = -876.3333
$ws.Range("H94").Value = 5435.136
$ws.Range("I94").Value = 6209.7646
$ws.Range("J94").Value = 2801.4
$ws.Range("K94").Value = 6209.7646
$ws.Range("L94").Value = 2801.4
$ws.Range("M94").Value = -5758.7646
$ws.Range("N94").Value = -3703.4
$ws.Range("H134").Value = 2052.8462
$ws.Range("I134").Value = 1585.4193
$ws.Range("J134").Value = 3864.125
$ws.Range("K134").Value = 4756.257900000001
$ws.Range("L134").Value = 11592.375
$ws.Range("M134").Value = -2221.257900000001
$ws.Range("N134").Value = -16662.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1677.591
$ws.Range("I16").Value = 1517.5385
$ws.Range("J16").Value = 1908.7778
$ws.Range("K16").Value = 1517.5385
$ws.Range("L16").Value = 1908.7778
$ws.Range("M16").Value = -1230.5385
$ws.Range("N16").Value = -2482.7778
$ws.Range("H31").Value = 10005801
$ws.Range("I31").Value = 33336066
$ws.Range("J31").Value = 7115.5713
$ws.Range("K31").Value = 33336066
$ws.Range("L31").Value = 7115.5713
$ws.Range("M31").Value = -33335771
$ws.Range("N31").Value = -7705.5713
$ws.Range("H34").Value = 10005801
$ws.Range("I34").Value = 33336066
$ws.Range("J34").Value = 7115.5713
$ws.Range("K34").Value = 33336066
$ws.Range("L34").Value = 7115.5713
$ws.Range("M34").Value = -33335864
$ws.Range("N34").Value = -7519.5713
$ws.Range("H58").Value = 13784.223
$ws.Range("I58").Value = 2053.7222
$ws.Range("J58").Value = 37245.223
$ws.Range("K58").Value = 2053.7222
$ws.Range("L58").Value = 37245.223
$ws.Range("M58").Value = -1850.7222
$ws.Range("N58").Value = -37651.223
$ws.Range("H105").Value = 1331.6666
$ws.Range("I105").Value = 1331.6666
$ws.Range("K105").Value = 1331.6666
$ws.Range("M105").Value = 415.3334
$ws.Range("H113").Value = 1677.591
$ws.Range("I113").Value = 1517.5385
$ws.Range("J113").Value = 1908.7778
$ws.Range("K113").Value = 1517.5385
$ws.Range("L113").Value = 1908.7778
$ws.Range("M113").Value = 652.4614999999999
$ws.Range("N113").Value = -6248.7778
$ws.Range("H136").Value = 13784.223
$ws.Range("I136").Value = 2053.7222
$ws.Range("J136").Value = 37245.223
$ws.Range("K136").Value = 6161.1666
$ws.Range("L136").Value = 111735.669
$ws.Range("M136").Value = -3611.1666
$ws.Range("N136").Value = -116835.669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3479.182
$ws.Range("I129").Value = 2947.25
$ws.Range("J129").Value = 3783.1428
$ws.Range("K129").Value = 8841.75
$ws.Range("L129").Value = 11349.4284
$ws.Range("M129").Value = -3841.75
$ws.Range("N129").Value = -21349.4284
$ws.Range("H131").Value = 1710.5834
$ws.Range("I131").Value = 1558
$ws.Range("J131").Value = 1802.1333
$ws.Range("K131").Value = 4674
$ws.Range("L131").Value = 5406.3999
$ws.Range("M131").Value = 366
$ws.Range("N131").Value = -15486.3999
$ws.Range("H133").Value = 6305.625
$ws.Range("J133").Value = 7871.2
$ws.Range("L133").Value = 23613.6
$ws.Range("N133").Value = -33733.6
$ws.Range("H137").Value = 3574.2727
$ws.Range("I137").Value = 4450
$ws.Range("K137").Value = 13350
$ws.Range("M137").Value = -8250

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2261.1
$ws.Range("J113").Value = 2376.375
$ws.Range("L113").Value = 2376.375
$ws.Range("N113").Value = -6716.375
$ws.Range("H132").Value = 2147.2188
$ws.Range("I132").Value = 2183.8965
$ws.Range("K132").Value = 6551.689499999999
$ws.Range("M132").Value = -4021.689499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2507
$ws.Range("I7").Value = 2042.6666
$ws.Range("K7").Value = 2042.6666
$ws.Range("M7").Value = -1930.6666
$ws.Range("H46").Value = 3964.2144
$ws.Range("I46").Value = 841.4286
$ws.Range("J46").Value = 7087
$ws.Range("K46").Value = 841.4286
$ws.Range("L46").Value = 7087
$ws.Range("M46").Value = -653.4286
$ws.Range("N46").Value = -7463
$ws.Range("H55").Value = 973.15
$ws.Range("I55").Value = 382.5
$ws.Range("J55").Value = 1563.8
$ws.Range("K55").Value = 382.5
$ws.Range("L55").Value = 1563.8
$ws.Range("M55").Value = -209.5
$ws.Range("N55").Value = -1909.8
$ws.Range("H122").Value = 6738.4
$ws.Range("J122").Value = 7430.125
$ws.Range("L122").Value = 22290.375
$ws.Range("N122").Value = -27190.375
$ws.Range("H126").Value = 2507
$ws.Range("I126").Value = 2042.6666
$ws.Range("K126").Value = 6127.9998
$ws.Range("M126").Value = -3657.9998
$ws.Range("H127").Value = 138999.25
$ws.Range("J127").Value = 138999.25
$ws.Range("L127").Value = 138999.25
$ws.Range("N127").Value = -148919.25
$ws.Range("H136").Value = 3774
$ws.Range("I136").Value = 3774
$ws.Range("K136").Value = 11322
$ws.Range("M136").Value = -8772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 47587
$ws.Range("J16").Value = 47587
$ws.Range("L16").Value = 47587
$ws.Range("N16").Value = -48171
$ws.Range("H46").Value = 136996
$ws.Range("J46").Value = 136996
$ws.Range("L46").Value = 136996
$ws.Range("N46").Value = -137458
$ws.Range("H62").Value = 4276.4
$ws.Range("I62").Value = 4132.6665
$ws.Range("J62").Value = 4492
$ws.Range("K62").Value = 4132.6665
$ws.Range("L62").Value = 4492
$ws.Range("M62").Value = -3508.6665
$ws.Range("N62").Value = -5740
$ws.Range("H65").Value = 4276.4
$ws.Range("I65").Value = 4132.6665
$ws.Range("J65").Value = 4492
$ws.Range("K65").Value = 20663.3325
$ws.Range("L65").Value = 22460
$ws.Range("M65").Value = -17543.3325
$ws.Range("N65").Value = -28700
$ws.Range("H132").Value = 1257.9395
$ws.Range("I132").Value = 880
$ws.Range("J132").Value = 5037.3335
$ws.Range("K132").Value = 2640
$ws.Range("L132").Value = 15112.0005
$ws.Range("M132").Value = -110
$ws.Range("N132").Value = -20172.0005
$ws.Range("H134").Value = 136996
$ws.Range("J134").Value = 136996
$ws.Range("L134").Value = 410988
$ws.Range("N134").Value = -416058
